# Updates crypto price/volume figures per the Sat Oct 19 19:45:20 UTC 2024 refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.198.18"
$ws.Range("E2").Value = "  -0.69%  "
$ws.Range("D3").Value = "2.644.50"
$ws.Range("E3").Value = "  -0.30%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "597.54"
$ws.Range("E5").Value = "  -0.37%  "
$ws.Range("D6").Value = "156.54"
$ws.Range("E6").Value = "  +1.13%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  -0.77%  "
$ws.Range("D9").Value = "'0.140"
$ws.Range("E9").Value = "  +2.09%  "
$ws.Range("E10").Value = "  -1.25%  "
$ws.Range("E11").Value = "  +0.60%  "
$ws.Range("D12").Value = "'0.350"
$ws.Range("E12").Value = "  +0.66%  "
$ws.Range("D13").Value = "27.98"
$ws.Range("E13").Value = "  +0.17%  "
$ws.Range("E14").Value = "  +0.56%  "
$ws.Range("D15").Value = "3.126.36"
$ws.Range("E15").Value = "  -0.32%  "
$ws.Range("D16").Value = "68.244.14"
$ws.Range("E16").Value = "  -0.43%  "
$ws.Range("D17").Value = "2.651.19"
$ws.Range("E17").Value = "  +0.04%  "
$ws.Range("D18").Value = "11.37"
$ws.Range("E18").Value = "  -0.43%  "
$ws.Range("D19").Value = "363.18"
$ws.Range("E19").Value = "  -1.00%  "
$ws.Range("E20").Value = "  -1.21%  "
$ws.Range("E21").Value = "  +3.31%  "
$ws.Range("D22").Value = "4.78"
$ws.Range("E22").Value = "  -1.74%  "
$ws.Range("E23").Value = "  -2.77%  "
$ws.Range("D24").Value = "75.45"
$ws.Range("E24").Value = "  +3.84%  "
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("E26").Value = "  -2.38%  "
$ws.Range("E28").Value = "  -1.16%  "
$ws.Range("E29").Value = "  +1.29%  "
$ws.Range("D30").Value = "556.41"
$ws.Range("E30").Value = "  -3.14%  "
$ws.Range("D31").Value = "8.06"
$ws.Range("E32").Value = "  -0.91%  "
$ws.Range("D33").Value = "1.85"
$ws.Range("E33").Value = "  -0.49%  "
$ws.Range("D34").Value = "0.999"
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("E35").Value = "  -1.92%  "
$ws.Range("D36").Value = "1.54"
$ws.Range("E36").Value = "  +0.20%  "
$ws.Range("D37").Value = "160.49"
$ws.Range("E37").Value = "  +1.12%  "
$ws.Range("D38").Value = "19.72"
$ws.Range("E38").Value = "  +2.47%  "
$ws.Range("E39").Value = "  +1.08%  "
$ws.Range("E40").Value = "  -3.38%  "
$ws.Range("D41").Value = "5.32"
$ws.Range("E41").Value = "  -1.22%  "
$ws.Range("E42").Value = "  +4.59%  "
$ws.Range("E43").Value = "  +0.29%  "
$ws.Range("E44").Value = "  -1.67%  "
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("D46").Value = "158.76"
$ws.Range("E46").Value = "  +1.39%  "
$ws.Range("E47").Value = "  -0.32%  "
$ws.Range("D48").Value = "21.96"
$ws.Range("E48").Value = "  +0.04%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "0.0782"
$ws.Range("E49").Value = "  +0.29%  "
$ws.Range("B50").Value = "Optimism"
$ws.Range("C50").Value = "https://coinranking.com/coin/n1p-s_gm1+optimism-op"
$ws.Range("D50").Value = "1.68"
$ws.Range("E50").Value = "  -1.96%  "
$ws.Range("E51").Value = "  -0.56%  "
